# Updated RAD Test Cases for Bay Restoration and Transportation for Month drop down.
# Re-run of the NewTaxReturn automated test cases: every row gets a fresh
# execution timestamp in column B (Date), and every previously-failing row
# now reports "Pass" in column A (Result).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewTaxReturn")

# Column A (Result): rows that previously reported "Fail" now report "Pass".
$passRows = @(18,19,21,22,24,25,26,27,28,29,30,31,33,34,36,37,38,39,40,41,42,43,45,46,48,49,50,51,52)
foreach ($r in $passRows) {
    $ws.Range("A$r").Value = "Pass"
}

# Column B (Date/CRN timestamp): every data row (2-52) re-stamped with the
# timestamp of the new test run.
$newDates = @(
    "Tue Mar 05 22:28:37 EST 2024",
    "Tue Mar 05 22:29:12 EST 2024",
    "Tue Mar 05 22:29:46 EST 2024",
    "Tue Mar 05 22:30:20 EST 2024",
    "Tue Mar 05 22:31:00 EST 2024",
    "Tue Mar 05 22:31:39 EST 2024",
    "Tue Mar 05 22:32:13 EST 2024",
    "Tue Mar 05 22:32:46 EST 2024",
    "Tue Mar 05 22:33:20 EST 2024",
    "Tue Mar 05 22:33:54 EST 2024",
    "Tue Mar 05 22:34:27 EST 2024",
    "Tue Mar 05 22:35:01 EST 2024",
    "Tue Mar 05 22:35:35 EST 2024",
    "Tue Mar 05 22:36:09 EST 2024",
    "Tue Mar 05 22:36:43 EST 2024",
    "Tue Mar 05 22:37:17 EST 2024",
    "Tue Mar 05 22:37:52 EST 2024",
    "Tue Mar 05 22:38:28 EST 2024",
    "Tue Mar 05 22:39:04 EST 2024",
    "Tue Mar 05 22:39:37 EST 2024",
    "Tue Mar 05 22:40:13 EST 2024",
    "Tue Mar 05 22:40:48 EST 2024",
    "Tue Mar 05 22:41:22 EST 2024",
    "Tue Mar 05 22:41:57 EST 2024",
    "Tue Mar 05 22:42:33 EST 2024",
    "Tue Mar 05 22:43:08 EST 2024",
    "Tue Mar 05 22:43:44 EST 2024",
    "Tue Mar 05 22:44:20 EST 2024",
    "Tue Mar 05 22:44:55 EST 2024",
    "Tue Mar 05 22:45:30 EST 2024",
    "Tue Mar 05 22:46:06 EST 2024",
    "Tue Mar 05 22:46:39 EST 2024",
    "Tue Mar 05 22:47:14 EST 2024",
    "Tue Mar 05 22:47:50 EST 2024",
    "Tue Mar 05 22:48:26 EST 2024",
    "Tue Mar 05 22:49:03 EST 2024",
    "Tue Mar 05 22:49:39 EST 2024",
    "Tue Mar 05 22:50:15 EST 2024",
    "Tue Mar 05 22:50:50 EST 2024",
    "Tue Mar 05 22:51:26 EST 2024",
    "Tue Mar 05 22:52:01 EST 2024",
    "Tue Mar 05 22:52:36 EST 2024",
    "Tue Mar 05 22:53:12 EST 2024",
    "Tue Mar 05 22:53:45 EST 2024",
    "Tue Mar 05 22:54:20 EST 2024",
    "Tue Mar 05 22:54:56 EST 2024",
    "Tue Mar 05 22:55:30 EST 2024",
    "Tue Mar 05 22:56:06 EST 2024",
    "Tue Mar 05 22:56:41 EST 2024",
    "Tue Mar 05 22:57:17 EST 2024",
    "Tue Mar 05 22:57:52 EST 2024"
)

$row = 2
foreach ($d in $newDates) {
    $ws.Range("B$row").Value = $d
    $row = $row + 1
}
